$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 06:22"

# Row 62: Kazajistan (updated figures)
$ws.Range("B62").Value = 2860
$ws.Range("C62").Value = 25
$ws.Range("D62").Value = 725
$ws.Range("E62").Value = 2110

# Row 166: Polinesia Francesa (updated figures)
$ws.Range("B166").Value = 58
$ws.Range("D166").Value = 49
$ws.Range("E166").Value = 9

# Row 168: Republica del Chad -> Republica de Africa Central
$ws.Range("A168").Value = "Republica de Africa Central"
$ws.Range("B168").Value = 50
$ws.Range("D168").Value = 10
$ws.Range("E168").Value = 40

# Row 169: Macao -> Republica del Chad
$ws.Range("A169").Value = "Republica del Chad"
$ws.Range("B169").Value = 46
$ws.Range("D169").Value = 15
$ws.Range("E169").Value = 31
$ws.Range("F169").Value = 0

# Row 170: Siria -> Macao
$ws.Range("A170").Value = "Macao"
$ws.Range("B170").Value = 45
$ws.Range("D170").Value = 32
$ws.Range("E170").Value = 13
$ws.Range("F170").Value = 1
$ws.Range("H170").Value = 0

# Row 171: Puerto Rico -> Siria
$ws.Range("A171").Value = "Siria"
$ws.Range("B171").Value = 43
$ws.Range("D171").Value = 19
$ws.Range("E171").Value = 21
$ws.Range("H171").Value = 3

# Row 172: Eritrea -> Puerto Rico
$ws.Range("A172").Value = "Puerto Rico"
$ws.Range("D172").Value = 1
$ws.Range("E172").Value = 36
$ws.Range("H172").Value = 2

# Row 173: Mongolia -> Eritrea
$ws.Range("A173").Value = "Eritrea"
$ws.Range("B173").Value = 39
$ws.Range("D173").Value = 13
$ws.Range("E173").Value = 26

# Row 174: San Martin (Parte Francesa) -> Mongolia
$ws.Range("A174").Value = "Mongolia"
$ws.Range("D174").Value = 10
$ws.Range("E174").Value = 28
$ws.Range("F174").Value = 0
$ws.Range("H174").Value = 0

# Row 175: Malaui -> San Martin (Parte Francesa)
$ws.Range("A175").Value = "San Martin (Parte Francesa)"
$ws.Range("B175").Value = 38
$ws.Range("D175").Value = 24
$ws.Range("E175").Value = 11
$ws.Range("F175").Value = 3

# Row 176: Guam -> Malaui
$ws.Range("A176").Value = "Malaui"
$ws.Range("B176").Value = 36
$ws.Range("D176").Value = 4
$ws.Range("E176").Value = 29
$ws.Range("F176").Value = 1
$ws.Range("H176").Value = 3

# Row 177: Zimbabue -> Guam
$ws.Range("A177").Value = "Guam"
$ws.Range("C177").Value = 0
$ws.Range("D177").Value = 0
$ws.Range("E177").Value = 31
$ws.Range("H177").Value = 1

# Row 178: Angola -> Zimbabue
$ws.Range("A178").Value = "Zimbabue"
$ws.Range("B178").Value = 32
$ws.Range("C178").Value = 1
$ws.Range("D178").Value = 5
$ws.Range("E178").Value = 23
$ws.Range("H178").Value = 4

# Row 179: Timor Oriental -> Angola
$ws.Range("A179").Value = "Angola"
$ws.Range("B179").Value = 27
$ws.Range("D179").Value = 6
$ws.Range("E179").Value = 19
$ws.Range("H179").Value = 2

# Row 180: Antigua y Barbuda -> Timor Oriental
$ws.Range("A180").Value = "Timor Oriental"
$ws.Range("D180").Value = 2
$ws.Range("E180").Value = 22
$ws.Range("F180").Value = 0
$ws.Range("H180").Value = 0

# Row 181: Botsuana -> Antigua y Barbuda
$ws.Range("A181").Value = "Antigua y Barbuda"
$ws.Range("B181").Value = 24
$ws.Range("D181").Value = 11
$ws.Range("E181").Value = 10
$ws.Range("F181").Value = 1
$ws.Range("H181").Value = 3

# Row 182: Laos -> Botsuana
$ws.Range("A182").Value = "Botsuana"
$ws.Range("B182").Value = 22
$ws.Range("D182").Value = 0
$ws.Range("E182").Value = 21
$ws.Range("H182").Value = 1

# Row 183: Republica de Africa Central -> Laos
$ws.Range("A183").Value = "Laos"
$ws.Range("D183").Value = 7
$ws.Range("E183").Value = 12
